# Remove footnote markers like " [1]".. " [5]" from vaccine/category names,
# and flatten embedded line breaks (two-line cells) into a single line with
# a space, across every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($n = 1; $n -le 5; $n++) {
        $token = "[" + $n + "]"
        [void]$ws.Cells.Replace($token, "")
    }
    [void]$ws.Cells.Replace("`n", " ")
}
